$p = $ppt.ActivePresentation

$oldText = "Are good employers good companies? " + [char]0x2013 + " Look within Top 100"
$newText = "Are good employers good companies? " + [char]0x2013 + " A look within Top 100"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldText) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}
